# Updated main GSC export data: the oldest day (2025-10-21) has aged out of
# the export window, so drop its row from the "Chart" data sheet. Excel
# shifts every subsequent row up by one, which is exactly what the new
# export reflects (last row 92 disappears, dimension shrinks to C91, and
# the now-unused "2025-10-21" shared string drops out on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 holds the 2025-10-21 data (Date / Non-HTTPS URLs / HTTPS URLs).
# Deleting it shifts rows 3..92 up into 2..91.
$ws.Rows.Item(2).Delete()
